# Links.xlsx update
# "more articles added with font size changes to fix some bugs"
#
# - A handful of new numeric "counter" cells are added next to existing
#   article rows (B20, A22, D24) and a brand new row is appended (row 30,
#   cell G30).
# - The hyperlink cell G22 gets a red highlight fill (this introduces a new
#   fill + cellXf in styles.xml, same as the source workbook diff).
# - The active selection/view is moved to column H, cell H7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New numeric cells for newly catalogued articles -----------------
$ws.Range("B20").Value = 15
$ws.Range("A22").Value = 17
$ws.Range("D24").Value = 19

# New row 30 at the bottom of the table
$ws.Range("G30").Value = 24

# --- Highlight G22 (existing hyperlink cell) with a red fill ----------
# RGB(192,0,0) == 0xC00000 -> Excel COM Color is 0x00BBGGRR, so for a pure
# red/green=0/blue=0 color the numeric value is simply the red component.
$ws.Range("G22").Interior.Color = 192

# --- Update the sheet view: scroll so column H is visible and select H7
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H7").Select()
